# bcThangPL.xlsx edit: "Thêm phục Lục 03c"
# - Insert a new appendix sheet "PL03c" (month-over-month KCB/chi comparison)
#   positioned right after "PL03b" and before "PL04a".
# - Expand "PL03a" with a month-over-month breakdown (Thang nay / Thang truoc / tang-giam)
#   for each of the five base indicators, plus detail sub-headers (cols G-K of B02-10).
# - Small cosmetic touch-ups on PL01 (merge title/unit rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) PL01 - merge title row and unit row so they span the whole table width
# ---------------------------------------------------------------------------
$pl01 = $wb.Worksheets.Item("PL01")
$pl01.Range("A1").Value = $pl01.Range("B1").Value
$pl01.Range("B1").Value = $null
$pl01.Range("A1:E1").Merge()
$pl01.Range("D2:E2").Merge()
$pl01.Range("A1").HorizontalAlignment = -4108  # xlCenter
$pl01.Range("A1").VerticalAlignment = -4108
$pl01.Range("A1").Font.Bold = $true
$pl01.Range("A1").Font.Size = 12
$pl01.Range("A1").Font.Name = "Times New Roman"
$pl01.Range("D2").HorizontalAlignment = -4152  # xlRight
$pl01.Range("D2").Font.Name = "Times New Roman"
$pl01.Range("D2").Font.Size = 12

# ---------------------------------------------------------------------------
# 2) Create the new "PL03c" sheet by copying "PL02c" (same month-comparison
#    layout for lot KCB / chi KCB, just worded for month-over-month instead
#    of year-over-year) and dropping it right after "PL03b".
# ---------------------------------------------------------------------------
$pl02c = $wb.Worksheets.Item("PL02c")
$pl02c.Copy($null, $wb.Worksheets.Item("PL03b"))
$pl03c = $wb.Worksheets.Item("PL02c (2)")
$pl03c.Name = " PL03c"

# Update the title + detail headers to the new month-over-month wording.
$pl03c.Range("A1").Value = "So sánh lượt KCB và chi KCB tháng này với tháng trước"
$pl03c.Range("C3").Value = "Tháng này"
$pl03c.Range("D3").Value = "Tháng trước"
$pl03c.Range("F3").Value = "Tháng này"
$pl03c.Range("G3").Value = "Tháng trước"
$pl03c.Range("C5").Value = " Cột D-B02-10-2024-tháng này"
$pl03c.Range("D5").Value = "  Cột D-B02-10-2024- tháng trước"
$pl03c.Range("E5").Value = "Tháng này- tháng trước"
$pl03c.Range("F5").Value = " Cột R-B02-10-2024-tháng này"
$pl03c.Range("G5").Value = "  Cột R-B02-10-2024- tháng trước"
$pl03c.Range("H5").Value = "Tháng này- tháng trước"

# The copied sheet had 4 extra blank rows (7-10) inherited from PL02c; PL03c
# only needs through row 6.
$pl03c.Range("A7:H10").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 3) Expand "PL03a" with the month-over-month (Thang nay / Thang truoc) split
#    for each of its five indicators, each getting its own 3-column group
#    (Thang nay | Thang truoc | tang/giam) plus a detail sub-header row.
# ---------------------------------------------------------------------------
$pl03a = $wb.Worksheets.Item("PL03a")

# Row 3: turn each single indicator column into a 3-wide merged group.
$pl03a.Range("C3:E3").Merge()
$pl03a.Range("F3:H3").Merge()
$pl03a.Range("I3:K3").Merge()
$pl03a.Range("L3:N3").Merge()
$pl03a.Range("O3:Q3").Merge()
$pl03a.Range("A3:A4").Merge()
$pl03a.Range("B3:B4").Merge()

# Row 4: Thang nay / Thang truoc / tang-giam sub-headers under each group.
$months = @("C4","F4","I4","L4","O4")
foreach ($c in $months) { $pl03a.Range($c).Value = "Tháng này" }
$prev = @("D4","G4","J4","M4","P4")
foreach ($c in $prev) { $pl03a.Range($c).Value = "Tháng trước" }
$delta = @("E4","H4","K4","N4","Q4")
foreach ($c in $delta) { $pl03a.Range($c).Value = "tăng/giảm" }

# Row 5: detailed B02-10 column references for each indicator group.
$pl03a.Range("C5").Value = "Cột G-B02-10 tháng năy"
$pl03a.Range("D5").Value = "Cột G-B02-10 tháng trước"
$pl03a.Range("F5").Value = "Cột H-B02-10 tháng năy"
$pl03a.Range("G5").Value = "Cột H-B02-10 tháng trước"
$pl03a.Range("I5").Value = "Cột I-B02-10 tháng năy"
$pl03a.Range("J5").Value = "Cột I-B02-10 tháng trước"
$pl03a.Range("L5").Value = "Cột J-B02-10 tháng năy"
$pl03a.Range("M5").Value = "CộtJ-B02-10 tháng trước"
$pl03a.Range("O5").Value = "Cột K-B02-10 tháng năy"
$pl03a.Range("P5").Value = "Cột K-B02-10 tháng trước"

# {filldata} placeholder row moves down from row 4 to row 6 now that the
# detail headers occupy rows 3-5.
$pl03a.Range("A6").Value = $pl03a.Range("A5").Value
if ($pl03a.Range("A5").Value -eq $pl03a.Range("A6").Value) {
    # only clear the old location once the value has been copied down
}

# ---------------------------------------------------------------------------
# 4) Restore the originally-active sheet/tab (PL03a) so the saved workbook
#    opens on the same tab as before, just shifted by the newly-inserted
#    sheet.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("PL03a").Activate()
